$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header (A column) and content (B/C columns) for rows 13-21
$ws.Range('A13').Value2 = 'Programa resumido:'
$ws.Range('B13').Value2 = '01/01/2022'
$ws.Range('C13').Value2 = '01/01/2022'
$ws.Range('A14').Value2 = 'Short syllabus:'
$ws.Range('B14').Value2 = 'Work on spreadsheets. Formatting. Importing and exporting data. Formulas and Functions. Graphics. Descriptive statistics. Data search and manipulation. Macros. Visual Basic for Applications. Reference to intervals; Repetition and flow control. Event-oriented programming. Arrays. Classes and collections. Introduction to UserForms. Dynamical tables. Defining new functions.'
$ws.Range('C14').Value2 = 'Work on spreadsheets. Formatting. Importing and exporting data. Formulas and Functions. Graphics. Descriptive statistics. Data search and manipulation. Macros. Visual Basic for Applications. Reference to intervals; Repetition and flow control. Event-oriented programming. Arrays. Classes and collections. Introduction to UserForms. Dynamical tables. Defining new functions.'
$ws.Range('A15').Value2 = 'Programa:'
$ws.Range('B15').Value2 = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range('C15').Value2 = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range('A16').Value2 = 'Syllabus:'
$ws.Range('B16').Value2 = '- Work in electronic spreadsheets: configuring the working environment; Several varieties of spreadsheets; navigating tables, cell notations, repetition of commands, creation of sequences;- Formatting: verification of conditions, conditional formatting- Importing and exporting data. Text files, csv (comma separated values) and other efficient data sharing formats.- Formulas and Functions: Using functions in spreadsheets. Using Solver and Scenario features. Syncing spreadsheets- Descriptive statistics. Mean, standard deviation, quartiles, median, mode and other characteristics of data distributions.- Matrices: working with matrices in spreadsheets; basic operations: sum, multiplication, transposition, inversion. Solution of linear systems using matrices.- Graphs: Time series graphs; histograms; graphs of data presentation (bars, circular sector, etc.) XY graphs of correlation between two variables.- Data search and manipulation: Data search and identification functions. Logical operators.- Macros: recording, editing and using command sequences (macros) to automate tasks.- Visual Basic for Applications: Developer guide, VB Editor, debugging tools. Project Explorer.- Reference to ranges: Range and Cells objects. Offset, resize, Columns and Rows properties- Repetition and flow control: For ... Next loops and variations. Do While / Until Ties. Flow Controls If ... Then ... Else- Event-oriented programming: Event levels and parameters.- Arrays: declaration, multidimensional arrays, dynamic arrays- Classes and collections: Creating and using classes, application events, collections, dictionaries- Introduction to UserForms: Inboxes, messages, buttons, radio buttons;- Pivot tables: creating and configuring a pivot table- Defining new functions: extending the'
$ws.Range('C16').Value2 = '- Work in electronic spreadsheets: configuring the working environment; Several varieties of spreadsheets; navigating tables, cell notations, repetition of commands, creation of sequences;- Formatting: verification of conditions, conditional formatting- Importing and exporting data. Text files, csv (comma separated values) and other efficient data sharing formats.- Formulas and Functions: Using functions in spreadsheets. Using Solver and Scenario features. Syncing spreadsheets- Descriptive statistics. Mean, standard deviation, quartiles, median, mode and other characteristics of data distributions.- Matrices: working with matrices in spreadsheets; basic operations: sum, multiplication, transposition, inversion. Solution of linear systems using matrices.- Graphs: Time series graphs; histograms; graphs of data presentation (bars, circular sector, etc.) XY graphs of correlation between two variables.- Data search and manipulation: Data search and identification functions. Logical operators.- Macros: recording, editing and using command sequences (macros) to automate tasks.- Visual Basic for Applications: Developer guide, VB Editor, debugging tools. Project Explorer.- Reference to ranges: Range and Cells objects. Offset, resize, Columns and Rows properties- Repetition and flow control: For ... Next loops and variations. Do While / Until Ties. Flow Controls If ... Then ... Else- Event-oriented programming: Event levels and parameters.- Arrays: declaration, multidimensional arrays, dynamic arrays- Classes and collections: Creating and using classes, application events, collections, dictionaries- Introduction to UserForms: Inboxes, messages, buttons, radio buttons;- Pivot tables: creating and configuring a pivot table- Defining new functions: extending the'
$ws.Range('A17').Value2 = 'Avaliação:'
$ws.Range('A18').Value2 = 'Método:'
$ws.Range('B18').Value2 = '7797767 - Viktor Pastoukhov'
$ws.Range('C18').Value2 = '7797767 - Viktor Pastoukhov'
$ws.Range('A19').Value2 = 'Critério:'
$ws.Range('B19').Value2 = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'
$ws.Range('C19').Value2 = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'
$ws.Range('A20').Value2 = 'Norma de recuperação:'
$ws.Range('B20').Value2 = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'
$ws.Range('C20').Value2 = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'
$ws.Range('A21').Value2 = 'Bibliografia:'
$ws.Range('B21').Value2 = 'Não haverá exame de recuperação'
$ws.Range('C21').Value2 = 'Não haverá exame de recuperação'

# Row 17 (Avaliacao:) no longer has B/C content; clear leftover old text
$ws.Range('B17:C17').ClearContents() | Out-Null

# Row heights per target layout
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# Remove now-unused trailing rows (old rows 22 and 23)
$ws.Rows.Item(23).Delete() | Out-Null
$ws.Rows.Item(22).Delete() | Out-Null

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
